$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for account 004404342 / ADSON / 50283.4 (spreadsheet row 4).
# This shifts every following row up by one.
$ws.Rows(4).Delete()

# The row that used to hold account 004397124 / MURYLO (now row 5 after the
# deletion above) gets its Saldo updated.
$ws.Cells.Item(5, 3).Value = 24085.3

# The row that used to hold account 004498637 / TIAGO / 4635.65 (now row 8
# after the deletion above) is replaced with account 004479965 / DIEGO / 12000.
# Force column A to stay plain text (it starts with a leading zero) and then
# drop the temporary number-format override so no stray style survives.
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "004479965"
$ws.Cells.Item(8, 1).ClearFormats()

$ws.Cells.Item(8, 2).Value = "DIEGO"
$ws.Cells.Item(8, 3).Value = 12000
